# This script fills in the missing box-score detail columns (D, F, G, I, J, K, L)
# for rows 51-58 on Sheet1, mirroring the pattern already present in rows 2-50:
#   D = Away Pts, F = Home Pts, G = Overtime flag ("NA"), H = Arena (already present)
#   I = Winning team, J = Losing team, K = Forecasted winner, L = Correct? (Yes/No)
# When the forecast (K) matches the actual winner (I), K gets a green highlight fill,
# matching the existing convention used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctFillColor = 5287936  # RGB(0,176,80) == theme fill used for correct picks (matches FF00B050)

$rows = @(
    @{ Row = 51; D = 148; F = 143; G = "NA"; I = "Dallas Mavericks";     J = "Atlanta Hawks";      K = "Atlanta Hawks";      L = "No";  Correct = $false },
    @{ Row = 52; D = 138; F = 104; G = "NA"; I = "Houston Rockets";      J = "Charlotte Hornets";  K = "Charlotte Hornets";  L = "No";  Correct = $false },
    @{ Row = 53; D = 131; F = 133; G = "NA"; I = "Indiana Pacers";       J = "Phoenix Suns";       K = "Indiana Pacers";     L = "Yes"; Correct = $true  },
    @{ Row = 54; D = 127; F = 107; G = "NA"; I = "Los Angeles Clippers"; J = "Toronto Raptors";    K = "Toronto Raptors";    L = "No";  Correct = $false },
    @{ Row = 55; D = 106; F = 107; G = "NA"; I = "Memphis Grizzlies";    J = "Orlando Magic";      K = "Memphis Grizzlies";  L = "Yes"; Correct = $true  },
    @{ Row = 56; D = 112; F = 100; G = "NA"; I = "Cleveland Cavaliers";  J = "Milwaukee Bucks";    K = "Cleveland Cavaliers"; L = "Yes"; Correct = $true  },
    @{ Row = 57; D = 107; F = 83;  G = "NA"; I = "Oklahoma City Thunder"; J = "New Orleans Pelicans"; K = "New Orleans Pelicans"; L = "No"; Correct = $false },
    @{ Row = 58; D = 100; F = 116; G = "NA"; I = "San Antonio Spurs";    J = "Portland Trail Blazers"; K = "San Antonio Spurs"; L = "Yes"; Correct = $true  }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    if ($r.Correct) {
        $ws.Range("K$n").Interior.Color = $correctFillColor
    }
    $ws.Range("L$n").Value = $r.L
}
